# Auto-generated edit script: update "想去人数" (F column) counts
# across sheets "展览" (Exhibition), "演出" (Performance), "全部类型" (All types)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4950  # was 4945
$ws.Range("F3").Value = 2773  # was 2771
$ws.Range("F5").Value = 2842  # was 2838
$ws.Range("F8").Value = 494  # was 492
$ws.Range("F9").Value = 1747  # was 1748
$ws.Range("F10").Value = 761  # was 758
$ws.Range("F11").Value = 500  # was 499
$ws.Range("F13").Value = 438  # was 435
$ws.Range("F14").Value = 1082  # was 1081
$ws.Range("F16").Value = 12  # was 10
$ws.Range("F19").Value = 1051  # was 1048
$ws.Range("F22").Value = 683  # was 678
$ws.Range("F23").Value = 759  # was 758
$ws.Range("F24").Value = 157  # was 155
$ws.Range("F25").Value = 15  # was 13
$ws.Range("F27").Value = 568  # was 561
$ws.Range("F28").Value = 64  # was 59
$ws.Range("F29").Value = 1668  # was 1666
$ws.Range("F30").Value = 1698  # was 1672
$ws.Range("F31").Value = 418  # was 412
$ws.Range("F32").Value = 47  # was 46
$ws.Range("F33").Value = 1598  # was 1587
$ws.Range("F34").Value = 230  # was 226
$ws.Range("F35").Value = 2423  # was 2410
$ws.Range("F36").Value = 420  # was 417
$ws.Range("F37").Value = 30  # was 29
$ws.Range("F38").Value = 631  # was 629
$ws.Range("F39").Value = 120  # was 119
$ws.Range("F40").Value = 72  # was 71
$ws.Range("F42").Value = 828  # was 823
$ws.Range("F43").Value = 1529  # was 1519
$ws.Range("F44").Value = 239  # was 235
$ws.Range("F47").Value = 76  # was 71
$ws.Range("F49").Value = 122  # was 121

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 108  # was 106

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4950  # was 4945
$ws.Range("F3").Value = 2773  # was 2771
$ws.Range("F4").Value = 2842  # was 2838
$ws.Range("F5").Value = 1747  # was 1748
$ws.Range("F8").Value = 761  # was 758
$ws.Range("F9").Value = 500  # was 499
$ws.Range("F11").Value = 438  # was 435
$ws.Range("F12").Value = 1082  # was 1081
$ws.Range("F15").Value = 1051  # was 1048
$ws.Range("F17").Value = 683  # was 678
$ws.Range("F18").Value = 759  # was 758
$ws.Range("F19").Value = 157  # was 155
$ws.Range("F20").Value = 108  # was 106
$ws.Range("F21").Value = 108  # was 106
$ws.Range("F23").Value = 15  # was 13
$ws.Range("F26").Value = 568  # was 561
$ws.Range("F27").Value = 1668  # was 1666
$ws.Range("F28").Value = 1698  # was 1672
$ws.Range("F29").Value = 418  # was 412
$ws.Range("F30").Value = 47  # was 46
$ws.Range("F33").Value = 2423  # was 2410
$ws.Range("F34").Value = 420  # was 417
$ws.Range("F38").Value = 30  # was 29
$ws.Range("F40").Value = 120  # was 119
$ws.Range("F41").Value = 72  # was 71
$ws.Range("F43").Value = 828  # was 823
$ws.Range("F44").Value = 1529  # was 1519
$ws.Range("F46").Value = 239  # was 235
$ws.Range("F48").Value = 76  # was 71

